# K means clustering is applied
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4914
$ws.Range("C2").Value = 0.4914
$ws.Range("D2").Value = 0.5338000059127808
$ws.Range("E2").Value = 0.9686999917030334
$ws.Range("F2").Value = 0.1958000063896179

$ws.Range("B3").Value = 0.5916
$ws.Range("C3").Value = 0.5916
$ws.Range("D3").Value = 0.6016
$ws.Range("E3").Value = 0.8482999801635742
$ws.Range("F3").Value = 0.2239000052213669

$ws.Range("B4").Value = 0.5077
$ws.Range("C4").Value = 0.5077
$ws.Range("D4").Value = 0.508
$ws.Range("E4").Value = 0.7921000123023987
$ws.Range("F4").Value = 0.2542999982833862

$ws.Range("B5").Value = 0.3642
$ws.Range("C5").Value = 0.3642
$ws.Range("D5").Value = 0.361
$ws.Range("E5").Value = 0.7243000268936157
$ws.Range("F5").Value = 0.1677999943494797

$ws.Range("B6").Value = 0.3564
$ws.Range("C6").Value = 0.3564
$ws.Range("D6").Value = 0.3484
$ws.Range("E6").Value = 0.7688000202178955
$ws.Range("F6").Value = 0.2064000070095062

$ws.Range("B7").Value = 0.3497
$ws.Range("C7").Value = 0.3497
$ws.Range("D7").Value = 0.3607000112533569
$ws.Range("E7").Value = 0.5302000045776367
$ws.Range("F7").Value = 0.2302999943494797

$ws.Range("B8").Value = 0.2914
$ws.Range("C8").Value = 0.2914
$ws.Range("D8").Value = 0.28
$ws.Range("E8").Value = 0.4905000030994415
$ws.Range("F8").Value = 0.214599996805191
